# Update to framework 7.0M2 - add BackTesting formula/attribute/label rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Attributes Labels" sheet: new row for the LastUpdateBackTesting label
# ---------------------------------------------------------------------------
$wsLabels = $wb.Worksheets.Item("Attributes Labels")
$wsLabels.Range("A3:D3").Copy($wsLabels.Range("A4:D4")) | Out-Null

$wsLabels.Range("A4").Value = "Properties_Predictive"
$wsLabels.Range("B4").Value = "LastUpdateBackTesting"
$wsLabels.Range("C4").Value = "LastUpdateBackTesting"
$wsLabels.Range("D4").Value = "LastUpdateBackTesting"

# ---------------------------------------------------------------------------
# 2) "Attributes" sheet: new row for the LastUpdateBackTesting attribute
# ---------------------------------------------------------------------------
$wsAttr = $wb.Worksheets.Item("Attributes")
$wsAttr.Range("A3:AE3").Copy($wsAttr.Range("A4:AE4")) | Out-Null

$wsAttr.Range("A4").Value = "Properties_Predictive"
$wsAttr.Range("B4").Value = "LastUpdateBackTesting"
$wsAttr.Range("C4").Value = "Date"
$wsAttr.Range("D4").Value = "Early Warning Properties"
$wsAttr.Range("E4").Value = 6
$wsAttr.Range("F4").Value = "N"
$wsAttr.Range("G4").Value = "N"
$wsAttr.Range("H4").Value = "N"
$wsAttr.Range("I4").Value = "Y"
$wsAttr.Range("J4").Value = "AUTO"
$wsAttr.Range("O4").Value = "Y"
$wsAttr.Range("Q4").Value = "UpdateBackTesting"
$wsAttr.Range("R4").Value = "N"
$wsAttr.Range("S4").Value = "N"
$wsAttr.Range("T4").Value = "Hourly"
$wsAttr.Range("Y4").Value = "None"
$wsAttr.Range("AE4").Value = "Y"

$wsAttr.Activate() | Out-Null
$wsAttr.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) "Formula" sheet: new row (UpdateBackTesting / RETURN LIB_RISK_NEW.backTesting();)
# ---------------------------------------------------------------------------
$wsFormula = $wb.Worksheets.Item("Formula")
$wsFormula.Range("A3:B3").Copy($wsFormula.Range("A4:B4")) | Out-Null
$wsFormula.Range("A4").Value = "UpdateBackTesting"
$wsFormula.Range("B4").Value = "RETURN LIB_RISK_NEW.backTesting();"

$wsFormula.Activate() | Out-Null
$wsFormula.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore "Attributes Labels" as the active sheet/tab and set its new
# selection, matching the final UI state from the diff.
# ---------------------------------------------------------------------------
$wsLabels.Activate() | Out-Null
$wsLabels.Range("B5").Select() | Out-Null
